$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A45").Value = "InRange"
$ws.Range("B45").Value = "Values"
$ws.Range("C45").Value = "Has Value"
$ws.Range("D45").Value = "Result"

$ws.Range("B46").Value = 1
$ws.Range("C46").Value = 1
$ws.Range("D46").Formula = "=COUNTIF(B46:B48,C46)"

$ws.Range("B47").Value = 2
$ws.Range("B48").Value = 3

$ws.Range("D46").Select() | Out-Null

